$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @"
0.04497535154223442,0.9920685887336731,0.2917045950889587,0.8868468403816223
0.007177175488322973,0.9987097382545471,0.2001671195030212,0.921319305896759
0.003940891940146685,0.9989563822746277,0.1297699809074402,0.9584740996360779
0.002107171574607491,0.9993928074836731,0.09903871268033981,0.9697993397712708
0.001605705823749304,0.9995825290679932,0.04621319100260735,0.9913570284843445
0.0009117427398450673,0.999772310256958,0.01477428246289492,0.9968209862709045
0.0009977055015042424,0.9997533559799194,0.03499288111925125,0.9891714453697205
0.0008631933014839888,0.9998292326927185,0.09147689491510391,0.9661235809326172
0.0005225273198448122,0.9998102784156799,0.09271353483200073,0.9603615999221802
0.001137578277848661,0.9997153878211975,0.1710263341665268,0.9370157122612
0.0007787683280184865,0.9997912645339966,0.2317436188459396,0.9198291301727295
0.0002849797310773283,0.9999810457229614,0.4951004385948181,0.8485992550849915
0.0004299433203414083,0.999886155128479,0.06532550603151321,0.9806278347969055
0.0008637411519885063,0.9997533559799194,0.05508442223072052,0.9809259176254272
0.000183705851668492,0.9999241232872009,0.06350543349981308,0.9778462052345276
0.0004025873786304146,0.9999051094055176,0.008791838772594929,0.9977150559425354
0.0001820193720050156,0.9999241232872009,0.06296680122613907,0.9780449271202087
0.0001363454648526385,0.9999620318412781,0.006865153089165688,0.9981124401092529
0.0006020630826242268,0.9999051094055176,0.3233460783958435,0.9206238985061646
0.0002294978476129472,0.999886155128479,0.03642316162586212,0.9914563894271851
0.0007004347280599177,0.999886155128479,0.01608124747872353,0.9946354031562805
0.0001461661304347217,0.9999620318412781,0.02759123779833317,0.9910590052604675
0.001204341184347868,0.9998481869697571,0.1449002623558044,0.9336379766464233
8.308889664476737e-05,0.9999620318412781,0.07274077832698822,0.9635406136512756
0.0004667158937081695,0.999886155128479,0.01894507184624672,0.9929465651512146
2.752896398305893e-05,1,0.02721733227372169,0.9885753989219666
1.169698498415528e-05,1,0.02653601951897144,0.9894695281982422
0.0001365855132462457,0.9999241232872009,0.05210800841450691,0.9851977229118347
0.0002549310156609863,0.999886155128479,0.01493685413151979,0.995628833770752
0.0007826474611647427,0.9998481869697571,0.06970255821943283,0.9774488210678101
0.000284696783637628,0.9999620318412781,0.03373663872480392,0.9880786538124084
2.759637209237553e-05,0.9999810457229614,0.05513083934783936,0.9830121397972107
0.000474411528557539,0.9999241232872009,0.2142944633960724,0.9417842030525208
0.0004598958184942603,0.9999620318412781,0.3564826250076294,0.9183389544487
0.0001915312459459528,0.9999430775642395,0.2683308124542236,0.9412875175476074
0.0001399757748004049,0.9999430775642395,0.1052131280303001,0.9665209650993347
0.00010727238986874,0.9999620318412781,0.2032349109649658,0.9440691471099854
0.000206468248507008,0.9998672008514404,0.1766960471868515,0.9531095027923584
0.0003268738510087132,0.9999620318412781,0.1182726249098778,0.9667196273803711
0.0005490024341270328,0.9999051094055176,0.1321217864751816,0.9637393355369568
0.000270524004008621,0.9999241232872009,0.09445898979902267,0.9731770157814026
7.746432675048709e-05,0.9999810457229614,0.7963477969169617,0.8204848170280457
0.0002391024609096348,0.9999051094055176,0.7631086707115173,0.8805881142616272
0.000231511308811605,0.9999810457229614,0.2004449516534805,0.9561891555786133
1.142250857810723e-05,1,0.5915632843971252,0.8705543279647827
0.0002372705639572814,0.9999430775642395,0.3228183686733246,0.9372143745422363
0.0002884367131628096,0.9999051094055176,0.03677201271057129,0.9926485419273376
8.320244887727313e-06,1,0.04546402394771576,0.9905622601509094
7.25709178368561e-05,0.9999810457229614,0.01960747502744198,0.9943373799324036
0.0004894349258393049,0.9999620318412781,0.02386931516230106,0.9922511577606201
"@

$rows = $data -split "`n" | Where-Object { $_.Trim() -ne "" }

$values = New-Object 'object[,]' 50,4
for ($i = 0; $i -lt $rows.Count; $i++) {
    $cols = $rows[$i].Trim() -split ","
    for ($j = 0; $j -lt 4; $j++) {
        $values[$i, $j] = [double]$cols[$j]
    }
}

$range = $ws.Range("A2:D51")
$range.Value = $values
